# Generate Report for Handback
#
# A new handback run produced a fresh file (fea4bf5d-bd52-4b86-bf17-c64116759e7a.md)
# with a new content hash (7c7ed79adc20636084129309279d5cb9fb2da81c) and new
# timestamps, replacing the previous 12ca941b-... run. At the same time the
# second tracked file was renamed from d19dcf06-dbd6-4f9c-bb07-c47f72314a55.md
# to ffffe74546cc-0a87-487a-ba6e-ba3bcbbca7de.md and now shares the same
# handback artifact (zh-cn/de-de xliff) as the first file. Update every sheet
# (Overview, zh-cn, de-de) accordingly, including the cell text and the
# display text cached on each hyperlink.

$wb = $excel.ActiveWorkbook

function Update-SheetHyperlinks {
    param(
        $ws,
        [hashtable]$displays   # A1 ref -> new display text (only for refs that changed)
    )

    # The underlying object model only lets us append hyperlinks, never patch
    # one in place, so capture every existing hyperlink's range + target
    # first, blow away the collection, then rebuild it — carrying forward the
    # original target address for untouched cells and the new display text
    # for the ones that changed.
    $refs = @()
    $addrs = @()
    $subs = @()
    $tips = @()
    $olddisplays = @()
    foreach ($h in $ws.Hyperlinks) {
        $refs += $h.Range.Address()
        $addrs += $h.Address
        $subs += $h.SubAddress
        $tips += $h.ScreenTip
        $olddisplays += $h.TextToDisplay
    }

    $ws.Hyperlinks.Delete()

    for ($i = 0; $i -lt $refs.Count; $i++) {
        $ref = $refs[$i]
        $disp = $olddisplays[$i]
        if ($displays.ContainsKey($ref)) {
            $disp = $displays[$ref]
        }
        $ws.Hyperlinks.Add($ws.Range($ref), $addrs[$i], $subs[$i], $tips[$i], $disp)
    }
}

# ---- Overview sheet ----------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.md"
$ws.Range("B2").Value = "e2e\fea4bf5d-bd52-4b86-bf17-c64116759e7a.md"
$ws.Range("G2").Value = "2016-09-01 19:13:24"

$ws.Range("A3").Value = "ffffe74546cc-0a87-487a-ba6e-ba3bcbbca7de.md"
$ws.Range("B3").Value = "e2e\ffffe74546cc-0a87-487a-ba6e-ba3bcbbca7de.md"
$ws.Range("G3").Value = "2016-09-01 19:13:24"

Update-SheetHyperlinks $ws @{
    '$B$2' = "e2e\fea4bf5d-bd52-4b86-bf17-c64116759e7a.md"
    '$B$3' = "e2e\ffffe74546cc-0a87-487a-ba6e-ba3bcbbca7de.md"
}

# ---- zh-cn sheet ---------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.md"
$ws.Range("G2").Value = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.7c7ed79adc20636084129309279d5cb9fb2da81c.zh-cn.xlf"
$ws.Range("H2").Value = "2016-09-01 19:13:19"
$ws.Range("I2").Value = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.md"
$ws.Range("J2").Value = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.7c7ed79adc20636084129309279d5cb9fb2da81c.zh-cn.xlf"
$ws.Range("K2").Value = "2016-09-01 19:13:38"

$ws.Range("A3").Value = "ffffe74546cc-0a87-487a-ba6e-ba3bcbbca7de.md"
$ws.Range("G3").Value = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.7c7ed79adc20636084129309279d5cb9fb2da81c.zh-cn.xlf"
$ws.Range("H3").Value = "2016-09-01 19:13:19"
$ws.Range("I3").Value = "ffffe74546cc-0a87-487a-ba6e-ba3bcbbca7de.md"
$ws.Range("J3").Value = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.7c7ed79adc20636084129309279d5cb9fb2da81c.zh-cn.xlf"
$ws.Range("K3").Value = "2016-09-01 19:13:38"

Update-SheetHyperlinks $ws @{
    '$A$2' = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.md"
    '$I$2' = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.md"
    '$A$3' = "ffffe74546cc-0a87-487a-ba6e-ba3bcbbca7de.md"
    '$I$3' = "ffffe74546cc-0a87-487a-ba6e-ba3bcbbca7de.md"
}

# ---- de-de sheet ---------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.md"
$ws.Range("G2").Value = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.7c7ed79adc20636084129309279d5cb9fb2da81c.de-de.xlf"
$ws.Range("H2").Value = "2016-09-01 19:13:24"
$ws.Range("I2").Value = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.md"
$ws.Range("J2").Value = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.7c7ed79adc20636084129309279d5cb9fb2da81c.de-de.xlf"
$ws.Range("K2").Value = "2016-09-01 19:13:46"

$ws.Range("A3").Value = "ffffe74546cc-0a87-487a-ba6e-ba3bcbbca7de.md"
$ws.Range("G3").Value = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.7c7ed79adc20636084129309279d5cb9fb2da81c.de-de.xlf"
$ws.Range("H3").Value = "2016-09-01 19:13:24"
$ws.Range("I3").Value = "ffffe74546cc-0a87-487a-ba6e-ba3bcbbca7de.md"
$ws.Range("J3").Value = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.7c7ed79adc20636084129309279d5cb9fb2da81c.de-de.xlf"
$ws.Range("K3").Value = "2016-09-01 19:13:46"

Update-SheetHyperlinks $ws @{
    '$A$2' = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.md"
    '$I$2' = "fea4bf5d-bd52-4b86-bf17-c64116759e7a.md"
    '$A$3' = "ffffe74546cc-0a87-487a-ba6e-ba3bcbbca7de.md"
    '$I$3' = "ffffe74546cc-0a87-487a-ba6e-ba3bcbbca7de.md"
}
